$d = $word.ActiveDocument

# Locate the paragraph containing "Changed home button layout to be scalable"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Changed home button layout to be scalable*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range

    # Insert a new bullet AFTER the target paragraph:
    # "Changed the home layout so the home box now is its own object"
    $rng.InsertAfter("`rChanged the home layout so the home box now is its own object")

    # Insert a new bullet BEFORE the target paragraph: "Added themed resumes"
    $rng.InsertBefore("Added themed resumes`r")
}

Write-Output "done"
